$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-ordered column headers in row 2: "StatusValuesID" moves to the front
# (column A) and "Notes" moves to the end (column L); every other column
# keeps its relative order and shifts accordingly.
$headers = @(
    "StatusValuesID",
    "ActivityBusinessKey",
    "DataVersionBusinessKey",
    "LocationBusinessKey",
    "OutcomeBusinessKey",
    "OutputBusinessKey",
    "ProgrammeBusinessKey",
    "ProjectBusinessKey",
    "ReportingPeriodBusinessKey",
    "StatusTypeBusinessKey",
    "SubOutputBusinessKey",
    "Notes"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $headers[$i]
}
